$wb = $excel.ActiveWorkbook

# 1. Remove the "Resumen_Posts" sheet entirely.
$resumen = $wb.Worksheets.Item("Resumen_Posts")
$resumen.Delete()

$ws = $wb.Worksheets.Item("Comentarios")

# 2a. H2:H10 date cells switch from date-only format to the same
#     date+time format used by column G (reuses style index 2).
$dtFormat = "YYYY-MM-DD HH:MM:SS"
foreach ($r in @(2,3,4,5,6,7,9,10)) {
    $ws.Cells.Item($r, 8).NumberFormat = $dtFormat
}

# 2b. J column (likes_count) switches from text-typed numbers to real numbers.
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(9, 10).Value = 2

# 2c. N8 raw scrape payload refreshed with newly re-scraped CDN URLs.
$n8New = @'
{'facebookUrl': 'https://www.facebook.com/reel/793063333529226', 'commentUrl': 'https://www.facebook.com/reel/793063333529226/?comment_id=790894940521711', 'id': 'Y29tbWVudDoxMjQ0MjA0NTc0NDE4NDQ2Xzc5MDg5NDk0MDUyMTcxMQ==', 'feedbackId': 'ZmVlZGJhY2s6MTI0NDIwNDU3NDQxODQ0Nl83OTA4OTQ5NDA1MjE3MTE=', 'date': '2025-10-13T15:02:06.000Z', 'attachments': [{'__typename': 'Sticker', 'animated_image': None, 'frame_count': 1, 'frame_rate': 83, 'frames_per_column': 1, 'frames_per_row': 1, 'label': "Avatar is smiling with their eyebrows upturned, they're resting their chin on their curled fists and little red hearts are floating around them.", 'pack': None, 'sprite_image': None, 'image': {'uri': 'https://scontent-ord5-2.xx.fbcdn.net/v/t39.1997-6/550754586_1301752704230824_6863359776300433879_n.webp?_nc_cat=105&ccb=1-7&_nc_sid=72b077&_nc_ohc=BwdGQZo9LcUQ7kNvwFeOvXh&_nc_oc=AdmN7vbaRQV2-zw7HeuI9TQL7EUJTiWktLw1V5MKCK7EoLNzT__i-8sadngpVJjeCVY&_nc_zt=26&_nc_ht=scontent-ord5-2.xx&_nc_gid=xv4cDyqxlfUsyGrfXq5gPw&oh=00_AfgVEFNGBF4RnkS-SEPiY4MjNyeYgZeXVYOn7gaiZJRUsw&oe=691458DB', 'width': 120, 'height': 120}, 'id': '2294760351042177'}], 'profileUrl': 'https://www.facebook.com/leidis.compartir', 'profilePicture': 'https://scontent-ord5-3.xx.fbcdn.net/v/t39.30808-1/436926704_122142271370189577_176725195454799216_n.jpg?stp=c7.0.212.212a_cp0_dst-jpg_s32x32_tt6&_nc_cat=107&ccb=1-7&_nc_sid=e99d92&_nc_ohc=RiZn7f5B2D0Q7kNvwHDSSbU&_nc_oc=AdlRkl61Z9QvD7f9GJulTUJ00vs1B3rVZ20CUNLJDLkD13vxaohIhdD37ArX1kPWUGA&_nc_zt=24&_nc_ht=scontent-ord5-3.xx&_nc_gid=xv4cDyqxlfUsyGrfXq5gPw&oh=00_AfgzHqBNwvB7PkOZoTb2ATCpAqxh7TfURZbpz-2uIiwLZA&oe=69144592', 'profileId': 'pfbid07TGb7jxGfMuSYSFug7CryiyzpwxVJSxgfyKovF9pez98BsUtq2aXCGwTpw7q97GNl', 'profileName': 'Leidis Compartir', 'likesCount': '0', 'threadingDepth': 0, 'facebookId': '1244204574418446', 'postTitle': 'Un solo cambio y el combo ahora sí está bueno, con Yogurt Alpina te va a encantar🙌', 'pageAdLibrary': {'id': '177587607187', 'woodhenge_creator_info': None}, 'inputUrl': 'https://www.facebook.com/reel/793063333529226'}
'@
$ws.Cells.Item(8, 14).Value = $n8New

# 2d. Three freshly scraped rows appended by the automation run.
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Instagram"
$ws.Cells.Item(12, 3).Value = "CAMPAÑA_MANUAL_MULTIPLE"
$ws.Cells.Item(12, 4).Value = "https://www.instagram.com/p/DPzKNF0DIqm/"
$ws.Cells.Item(12, 5).Value = "'"
$ws.Cells.Item(12, 6).Value = "'"
$ws.Cells.Item(12, 7).Value = "'"
$ws.Cells.Item(12, 8).Value = "'"
$ws.Cells.Item(12, 9).Value = "'"
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = $false
$ws.Cells.Item(12, 13).Value = "https://instagram.com/"
$n12 = @'
{'url': 'https://www.instagram.com/p/DPzKNF0DIqm/', 'requestErrorMessages': ['HTTP 200 undefined'], 'error': 'no_items', 'errorDescription': 'Empty or private data for provided input'}
'@
$ws.Cells.Item(12, 14).Value = $n12

$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Facebook"
$ws.Cells.Item(13, 3).Value = "CAMPAÑA_MANUAL_MULTIPLE"
$ws.Cells.Item(13, 4).Value = "https://www.facebook.com/reel/793063333529226"
$ws.Cells.Item(13, 5).Value = "'"
$ws.Cells.Item(13, 6).Value = "'"
$ws.Cells.Item(13, 7).Value = 45943.62645833333
$ws.Cells.Item(13, 7).NumberFormat = $dtFormat
$ws.Cells.Item(13, 8).Value = 45943
$ws.Cells.Item(13, 8).NumberFormat = $dtFormat
$ws.Cells.Item(13, 9).Value = "15:02:06"
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = $false
$ws.Cells.Item(13, 13).Value = "'"
$n13 = @'
{'facebookUrl': 'https://www.facebook.com/reel/793063333529226', 'commentUrl': 'https://www.facebook.com/reel/793063333529226/?comment_id=790894940521711', 'id': 'Y29tbWVudDoxMjQ0MjA0NTc0NDE4NDQ2Xzc5MDg5NDk0MDUyMTcxMQ==', 'feedbackId': 'ZmVlZGJhY2s6MTI0NDIwNDU3NDQxODQ0Nl83OTA4OTQ5NDA1MjE3MTE=', 'date': '2025-10-13T15:02:06.000Z', 'attachments': [{'__typename': 'Sticker', 'animated_image': None, 'frame_count': 1, 'frame_rate': 83, 'frames_per_column': 1, 'frames_per_row': 1, 'label': "Avatar is smiling with their eyebrows upturned, they're resting their chin on their curled fists and little red hearts are floating around them.", 'pack': None, 'sprite_image': None, 'image': {'uri': 'https://scontent-ord5-2.xx.fbcdn.net/v/t39.1997-6/550754586_1301752704230824_6863359776300433879_n.webp?_nc_cat=105&ccb=1-7&_nc_sid=72b077&_nc_ohc=BwdGQZo9LcUQ7kNvwGylSHl&_nc_oc=Adm2E9Dp3vkqWw45dHCMf3JWruj_yAKMcfTA9e5xl1TbNrgy5STSwO8Nt_7-vNescOo&_nc_zt=26&_nc_ht=scontent-ord5-2.xx&_nc_gid=kl03MC1Xy9Hh6m8EiUqcAA&oh=00_Afi6-VLryKvoJRQF8BW64QWEAtHnSY7nlk5xe9LdoazbYw&oe=691458DB', 'width': 120, 'height': 120}, 'id': '2294760351042177'}], 'profileUrl': 'https://www.facebook.com/leidis.compartir', 'profilePicture': 'https://scontent-ord5-3.xx.fbcdn.net/v/t39.30808-1/436926704_122142271370189577_176725195454799216_n.jpg?stp=c7.0.212.212a_cp0_dst-jpg_s32x32_tt6&_nc_cat=107&ccb=1-7&_nc_sid=e99d92&_nc_ohc=RiZn7f5B2D0Q7kNvwFaEddo&_nc_oc=AdmHwcOXKl27Qd6uEKlOa4WaLVKxcf0SZAFAiwWFtGv0z8DtHTRE0ChrvWxmC0X8aws&_nc_zt=24&_nc_ht=scontent-ord5-3.xx&_nc_gid=kl03MC1Xy9Hh6m8EiUqcAA&oh=00_AfifH-IdIG4AIe3mjMtPpm92Fusj9s12A84vUrB9KJZChw&oe=69144592', 'profileId': 'pfbid07TGb7jxGfMuSYSFug7CryiyzpwxVJSxgfyKovF9pez98BsUtq2aXCGwTpw7q97GNl', 'profileName': 'Leidis Compartir', 'likesCount': '0', 'threadingDepth': 0, 'facebookId': '1244204574418446', 'postTitle': 'Un solo cambio y el combo ahora sí está bueno, con Yogurt Alpina te va a encantar🙌', 'pageAdLibrary': {'id': '177587607187', 'woodhenge_creator_info': None}, 'inputUrl': 'https://www.facebook.com/reel/793063333529226'}
'@
$ws.Cells.Item(13, 14).Value = $n13

$ws.Cells.Item(14, 1).Value = 3
$ws.Cells.Item(14, 2).Value = "Instagram"
$ws.Cells.Item(14, 3).Value = "CAMPAÑA_MANUAL_MULTIPLE"
$ws.Cells.Item(14, 4).Value = "https://www.instagram.com/p/DPpVC6UjLlp/"
$ws.Cells.Item(14, 5).Value = "'"
$ws.Cells.Item(14, 6).Value = "'"
$ws.Cells.Item(14, 7).Value = "'"
$ws.Cells.Item(14, 8).Value = "'"
$ws.Cells.Item(14, 9).Value = "'"
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = $false
$ws.Cells.Item(14, 13).Value = "https://instagram.com/"
$n14 = @'
{'url': 'https://www.instagram.com/p/DPpVC6UjLlp/', 'requestErrorMessages': ['HTTP 200 undefined'], 'error': 'no_items', 'errorDescription': 'Empty or private data for provided input'}
'@
$ws.Cells.Item(14, 14).Value = $n14

Write-Output "edit applied"
